$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.963.14'
$ws.Range('E2').Value = '  -3.99%  '
$ws.Range('D3').Value = '1.869.57'
$ws.Range('E3').Value = '  -3.21%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('E5').Value = '  -2.62%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('D7').Value = '0.4353'
$ws.Range('E7').Value = '  -5.66%  '
$ws.Range('D8').Value = '0.3756'
$ws.Range('E8').Value = '  -2.06%  '
$ws.Range('D9').Value = '0.07476'
$ws.Range('E9').Value = '  -3.83%  '
$ws.Range('D10').Value = '0.9369'
$ws.Range('E10').Value = '  -4.51%  '
$ws.Range('D11').Value = '21.25'
$ws.Range('E11').Value = '  -5.85%  '
$ws.Range('D12').Value = '1.923.24'
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('D13').Value = '6.742'
$ws.Range('E13').Value = '  -3.60%  '
$ws.Range('E14').Value = '  -4.61%  '
$ws.Range('D15').Value = '0.06864'
$ws.Range('E15').Value = '  -2.94%  '
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '81.47'
$ws.Range('E17').Value = '  -3.68%  '
$ws.Range('D18').Value = '0.000009051'
$ws.Range('E18').Value = '  -5.31%  '
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').Value = '15.82'
$ws.Range('E20').Value = '  -5.79%  '
$ws.Range('D21').Value = '27.955.14'
$ws.Range('E21').Value = '  -4.00%  '
$ws.Range('D22').Value = '5.125'
$ws.Range('E22').Value = '  -4.28%  '
$ws.Range('D23').Value = '11.04'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').Value = '2.132.86'
$ws.Range('E24').Value = '  -1.65%  '
$ws.Range('D25').Value = '2.038'
$ws.Range('E25').Value = '  -2.12%  '
$ws.Range('D26').Value = '152.94'
$ws.Range('E26').Value = '  -3.31%  '
$ws.Range('D27').Value = '18.60'
$ws.Range('E27').Value = '  -2.81%  '
$ws.Range('D28').Value = '5.602'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('D29').Value = '113.42'
$ws.Range('E29').Value = '  -3.94%  '
$ws.Range('E30').Value = '  -8.24%  '
$ws.Range('D31').Value = '0.09026'
$ws.Range('E31').Value = '  -3.44%  '
$ws.Range('D32').Value = '0.8112'
$ws.Range('E32').Value = '  -6.43%  '
$ws.Range('D33').Value = '4.805'
$ws.Range('E33').Value = '  -6.64%  '
$ws.Range('E34').Value = '  -5.76%  '
$ws.Range('D35').Value = '2.972'
$ws.Range('E35').Value = '  -1.49%  '
$ws.Range('D36').Value = '1.002'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D38').Value = '1.116'
$ws.Range('E38').Value = '  -3.88%  '
$ws.Range('D39').Value = '0.01983'
$ws.Range('E39').Value = '  -3.34%  '
$ws.Range('D40').Value = '2.958'
$ws.Range('E40').Value = '  -3.43%  '
$ws.Range('D41').Value = '0.5270'
$ws.Range('E41').Value = '  -4.75%  '
$ws.Range('D42').Value = '0.1698'
$ws.Range('E42').Value = '  -3.40%  '
$ws.Range('D43').Value = '6.978'
$ws.Range('E43').Value = '  -7.71%  '
$ws.Range('D44').Value = '8.780'
$ws.Range('E44').Value = '  -6.55%  '
$ws.Range('D45').Value = '0.06752'
$ws.Range('E45').Value = '  -2.59%  '
$ws.Range('D46').Value = '0.4888'
$ws.Range('E46').Value = '  -6.32%  '
$ws.Range('D47').Value = '10.59'
$ws.Range('E47').Value = '  -5.63%  '
$ws.Range('D48').Value = '106.91'
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('D49').Value = '1.675'
$ws.Range('E49').Value = '  -6.12%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').Value = '1.001'
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '1.905'
$ws.Range('E51').Value = '  -14.01%  '
